# Edit: remove the "CEP: {{ cep }} " field block from the address paragraph
# and relocate the (auto-managed) "_GoBack" bookmark so that it now sits right
# after the "Bairro: {{ bairro }} " block instead of right before the
# "rua: {{ endereco }} " block.
#
# Net visible-text effect:
#   "... rua: {{ endereco }} Bairro: {{ bairro }} CEP: {{ cep }} Cidade: ..."
# becomes
#   "... rua: {{ endereco }} Bairro: {{ bairro }} Cidade: ..."

$d = $word.ActiveDocument

# --- Step 1: delete the "CEP: {{ cep }} " merge-field block -----------------
$found = $d.Content.Find.Execute(
    "CEP: {{ cep }} ",  # FindText
    $false,             # MatchCase
    $false,             # MatchWholeWord
    $false,             # MatchWildcards
    $false,             # MatchSoundsLike
    $false,             # MatchAllWordForms
    $true,              # Forward
    1,                  # Wrap (wdFindContinue)
    $false,             # Format
    "",                 # ReplaceWith (delete -> empty string)
    2                   # Replace (wdReplaceAll)
)

# --- Step 2: move the "_GoBack" bookmark to just after "{{ bairro }} " ------
$fullText = $d.Content.Text
$anchorText = "{{ bairro }} "
$idx = $fullText.IndexOf($anchorText)
if ($idx -ge 0) {
    $pos = $idx + $anchorText.Length
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
